# edit.ps1 - apply "Se aplica inyección de dependencias" changes to mejoras.docx
#
# Strategy: the Word OM exposed by this host doesn't give us a direct way to
# author w:proofErr markers or to split an existing run into several runs, so
# we drive everything through Range.InsertXML with a minimal OOXML "package"
# fragment. When the target Range spans real text, InsertXML replaces that
# text with the supplied run/paragraph markup; when the target Range is the
# zero-length point right before a paragraph mark, InsertXML inserts whole
# paragraphs in front of it (the final <w:p> in the fragment takes over the
# original paragraph mark). Both behaviours are used below.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($paragraph, [string]$innerParagraphXml) {
    $xml = New-PkgXml('<w:body>' + $innerParagraphXml + '</w:body>')
    $paragraph.Range.InsertXML($xml)
}

function Insert-ParagraphsBefore($paragraph, [string]$paragraphsXml) {
    $xml = New-PkgXml('<w:body>' + $paragraphsXml + '</w:body>')
    $start = $paragraph.Range.Start
    $pt = $d.Range($start, $start)
    $pt.InsertXML($xml)
}

# --- 1. "1. Cambiar datos estáticos a stream" -> split off "stream" ---------
$p1 = $d.Paragraphs(2)
$body1 = '<w:p>' +
    '<w:r><w:t>1.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Cambiar datos estáticos a </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>stream</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p1 $body1

# --- 2. "2. Cambiar log en texto plano a json (...funcionen los json...)" --
$p2 = $d.Paragraphs(6)
$body2 = '<w:p>' +
    '<w:r><w:t>2.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Cambiar log en texto plano a </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>json</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (he tenido que añadir un paquete extra para que funcione</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">n los </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>json</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, pero no ha habido mayores problemas. </w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p2 $body2

# --- 3. Replace the trailing empty paragraph with the new "3. Refactorizar
#        para aplicar inyección de dependencias" section -------------------
$lastParaIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastParaIndex)

$newParas =
    '<w:p><w:r><w:t>3. Refactorizar para aplicar inyección de dependencias</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>Ventajas:</w:t></w:r></w:p>' +
    '<w:p>' +
        '<w:r><w:t xml:space="preserve">- </w:t></w:r>' +
        '<w:r><w:t>Facilita las pruebas: Permite sustituir componentes reales por simulaciones (</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>mocks</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>) para probar la lógica de forma aislada.</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
        '<w:r><w:t xml:space="preserve">- </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">Aumenta la flexibilidad: Puedes cambiar una implementación (ej. </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>AlertSender</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> por </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>SMSAlertSender</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>) sin modificar las clases que la usan.</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
        '<w:r><w:t xml:space="preserve">- </w:t></w:r>' +
        '<w:r><w:t>Reduce el acoplamiento: Las clases son más independientes, lo que hace el código más fácil de entender, mantener y reutilizar.</w:t></w:r>' +
    '</w:p>'

Insert-ParagraphsBefore $pLast $newParas

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
